$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row before row 48, shifting rows 48:52 down to 49:53
$ws.Rows.Item(48).Insert()

# Fill the new row 48 with the new weekly record values
$ws.Cells.Item(48, 1).Value = 3
$ws.Cells.Item(48, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(48, 3).Value = "Coquimbo"
$ws.Cells.Item(48, 4).Value = 44753
$ws.Cells.Item(48, 5).Value = 5
$ws.Cells.Item(48, 6).Value = 100112035
$ws.Cells.Item(48, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(48, 8).Value = "Sin especificar"
$ws.Cells.Item(48, 9).Value = "Primera"
$ws.Cells.Item(48, 10).Value = 80
$ws.Cells.Item(48, 11).Value = 15000
$ws.Cells.Item(48, 12).Value = 16000
$ws.Cells.Item(48, 13).Value = 15500
$ws.Cells.Item(48, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(48, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(48, 16).Value = 1033
$ws.Cells.Item(48, 17).Value = 15
$ws.Cells.Item(48, 18).Value = "Hortaliza"

# Apply the same date number format as the other D-column cells (style index 2: YYYY-MM-DD HH:MM:SS)
$ws.Cells.Item(48, 4).NumberFormat = $ws.Cells.Item(49, 4).NumberFormat
